$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistic (C) and P-Value (D) figures per
# "Correcion a Diebold Mariano y revision de Cap1"

$ws.Range("C2").Value = -1.158246574716258
$ws.Range("D2").Value = 0.2591828707666679

$ws.Range("C3").Value = -0.9396386905507168
$ws.Range("D3").Value = 0.3576060494037518

$ws.Range("C4").Value = -0.3856587676128717
$ws.Range("D4").Value = 0.703453363008347

$ws.Range("C5").Value = 0.06174926957693898
$ws.Range("D5").Value = 0.9513200670531778

$ws.Range("C6").Value = 0.1637707048218131
$ws.Range("D6").Value = 0.8714068787166238

$ws.Range("C7").Value = 0.8116739843726634
$ws.Range("D7").Value = 0.4256740513852633

$ws.Range("C8").Value = 1.542299189524414
$ws.Range("D8").Value = 0.1372635383858023

$ws.Range("C9").Value = 0.6073900604025956
$ws.Range("D9").Value = 0.5498106075528155

$ws.Range("C10").Value = 0.7959716055095506
$ws.Range("D10").Value = 0.4345537940496706

$ws.Range("C11").Value = 0.4604424176370345
$ws.Range("D11").Value = 0.6497175147277927
